$d = $word.ActiveDocument

# --- Paragraph 1: title line + break + subtitle line ---
$p1 = $d.Paragraphs(1)
$full1 = $p1.Range
$txt1 = $full1.Text
$brPos = $txt1.IndexOf([char]11)
$start1 = $full1.Start

$titleRange = $d.Range($start1, $start1 + $brPos)
$titleRange.Text = "המאמר היומי של מייק:  14.04.25"

# re-fetch paragraph 1 since offsets shifted after the edit above
$p1b = $d.Paragraphs(1)
$full1b = $p1b.Range
$txt1b = $full1b.Text
$brPos2 = $txt1b.IndexOf([char]11)
$start1b = $full1b.Start
$subtitleRange = $d.Range($start1b + $brPos2 + 1, $full1b.End - 1)
$subtitleRange.Text = "Draft Model Knows When to Stop: A Self-Verification Length Policy for Speculative Decoding"

# --- Paragraph 2: intro ---
$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
$body2 = $d.Range($r2.Start, $r2.End - 1)
$body2.Text = "המאמר הזה משך את עיניי כבר בהסתכלות הראשונה בגלל צמד המילים ״Speculative Decoding״ או SD בקצרה שמאוד קרוב לליבי - אפילו הכנתי על זה מצגת די מקיפה שאני מציג אותה בפורומים שונים. SD מאפשר להגדיל את קצב גנרוט טקסט על ידי מודל שפה באמצעות שילוב מודל היעד עם מודל קטן מהיר יותר וכמובן יותר חלש ממודל היעד. המודל הקטן מייצר כמה טוקנים בצורה אוטורגרסיבית ומודל היעד חוזה מנצל טוקנים אלו כדי לחזות בו זמנית את הטוקנים הבאים שלו. זה מאפשר להגדיל את הקצב הדגימה של המודל הגדול בצורה ניכרת."

# --- Paragraph 3: context on conditioned diffusion / SD bottleneck ---
$p3 = $d.Paragraphs(3)
$r3 = $p3.Range
$body3 = $d.Range($r3.Start, $r3.End - 1)
$body3.Text = "השיטה מנצלת את העובדה שצוואר הבקבוק של תהליך הגנרוט העברת דאטה בין הזכרונות של gpu (בפרט HBM הגדול ואיטי ו-SRAM הקן אך מהיר בחלק החישובי של ה-gpu). אז SD מבצע חיזוי מהיר עם המודל הקטן ואז החיזוי הבו זמני על ידי המודל הגדול עם הטוקנים שנחזו על ידי המודל הקטן. אבל יש שם קאץ' כמובן: כדי לקבל את אותה התפלגות הטוקנים עם המודל הגדול דרך ניצול הטוקנים של המודל הקטן יש צורך בלבצע סוג של rejection sampling או RS."

# --- Paragraph 4: CFG goal / RS recap ---
$p4 = $d.Paragraphs(4)
$r4 = $p4.Range
$body4 = $d.Range($r4.Start, $r4.End - 1)
$body4.Text = "אזכיר ש-RS מאפשר לדגום מהתפלגות קלה לדגימה f כדי ליצור מדגם הדגום מהתפלגות אחרת g שקשה לדגום ממנה בצורה ישירה. אז אנו דוגמים נקודה x מ-f אז מקבלים את הדגימה בהסתברות השווה ליחס בין (f(x ל- (g(x (אם יחס זה גדול מ-1 הנקודה מתקבלת אוטומטית). ניתן להוכיח שנקודות שהנדגמות באופן זה מפולגת עם התפלגות הרצויה g."

# --- Paragraph 5: CFG mechanics / SD acceptance sampling ---
$p5 = $d.Paragraphs(5)
$r5 = $p5.Range
$body5 = $d.Range($r5.Start, $r5.End - 1)
$body5.Text = "אז במקרה שלנו (SD) אנו עושים משהו דומה עבור הטוקנים הנדגמים עם המודל קטן. במהלך השלב השני (דגימה בו זמנית מהמודל הגדול) עבור כל טוקן הנגדם מהמודל הקטן אנו מחשבים את היחס בין ההסתברויות של המודלים ואנו ״מקבלים״ את הטוקנים של המודל הקטן בהסתברות השווה ליחס סיבה. אחרי שהטוקן הראשון של המודל הקטן ״סורב״ (rejected) המודל הגדול מגנרט טוקן הבא עם המודל הגדול ואז המודל הקטן שוב מופעל לגנרט את הטוקנים הבאים. ד״א גם הטוקנים שמתקבלים מגונרטים עם עם התפלגות המחושבת משתי ההתפלגויות של הטוקן (של המודל הקטן ושל הגדול)."

# --- Paragraph 6: authors' finding / acceptance-rate control ---
$p6 = $d.Paragraphs(6)
$r6 = $p6.Range
$body6 = $d.Range($r6.Start, $r6.End - 1)
$body6.Text = "כמו שכבר הצלחתם להבין ״שליטה״ ב acceptance rate של טוקנים של המודל הקטן היא מאוד חשובה - באידיאל אנו רוצים לדגום מהמודל הקטן רק את הטוקנים שיתקבלו. המאמר מציע שיטה לשפר את ה-acceptance rate. המאמר מראה  שהממוצע של acceptance rate (זה די קל) שווה להפרש בין 1 למה שנקרא total variation distance או TBD בקצרה בין ההתפלגויות של שני המודלים (המותנות בהקשר). ולמזלנו עומד לרשותנו אי שוויון לא ידוע במיוחד שמאפשר לחסום TBD מלמטה עם הפרש בין קרוס-אנטרופי בין התפלגויותם של שני המודלים (עבור טוקן נתון בהינתן הקשרו) לאנטרופיה של טוקן של המודל הקטן."

# --- Paragraph 7: method summary (gains xml:space=preserve, trailing space) ---
$p7 = $d.Paragraphs(7)
$r7 = $p7.Range
$body7 = $d.Range($r7.Start, $r7.End - 1)
$body7.Text = "אבל כמובן שאנו לא יכולים לחשב את הקרוס אנטרופי בין התפלגויות אלו בשלב דגימה מהמודל הגדול עבור כל הטוקנים כי אנו דוגמים כל הטוקנים ממנו בו זמנית ולא יודעים מראש התפלגות מותנית של כל טוקן של המודל הגדול. אז המאמר ״משערך״ את הקרוס אנטרופי הזה על זמן מדגם די גדול דרך קבוע (קצת גדול מ 1) מוכפל באנטרופיה של הטוקן של המודל הקטן. אחרי שיש לנו את הקרוס-אנטרופי אנו יכולים לשערך את ה-acceptance rate עבור כל טוקן של המודל הקטן לפני הדגימה מהמודל הגדול. זה מאפשר לנו לקבוע את מספר הטוקנים מהמודל הקטן שעבורם תתבצע דגימה בו זמנית מהמודל הגדול - פשוט בוחרים טוקנים עד שה-acceptance rate המשוערך גבוה מאיזה סף. "

# --- Paragraph 8: closing remark ---
$p8 = $d.Paragraphs(8)
$r8 = $p8.Range
$body8 = $d.Range($r8.Start, $r8.End - 1)
$body8.Text = "רעיון נחמד אבל בחירת הקבוע בשלב האחרון לדעתי לא אופטימלית ואני מקווה שבקרוב ייצאו מחקרים המשפרים את ההיבט הזה של השיטה המוצעת."

# --- Paragraph 9: arxiv link ---
$p9 = $d.Paragraphs(9)
$r9 = $p9.Range
$body9 = $d.Range($r9.Start, $r9.End - 1)
$body9.Text = "https://arxiv.org/abs/2411.18462"
